$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; unprotect to make edits, then re-protect afterwards.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A40).
$ws.Range("A40").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-09 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for holdings rows 2-37.
$ws.Range("D2").Value = 0.03120144073608346
$ws.Range("E2").Value = 0.01357210179076351
$ws.Range("D3").Value = 0.03527780163059141
$ws.Range("E3").Value = 0.008402688860435381
$ws.Range("D4").Value = 0.03333525462845877
$ws.Range("E4").Value = -0.00707858428314323
$ws.Range("D5").Value = 0.03880976943282193
$ws.Range("E5").Value = 0.02209559603552247
$ws.Range("D6").Value = 0.01618004966351849
$ws.Range("E6").Value = -0.00625227190112676
$ws.Range("D7").Value = 0.01650823823902715
$ws.Range("E7").Value = 0.001638877012968276
$ws.Range("D8").Value = 0.03194462905868694
$ws.Range("E8").Value = 0.008384640311082325
$ws.Range("D9").Value = 0.03250125570359984
$ws.Range("E9").Value = 0.001954397394136942
$ws.Range("D10").Value = 0.03261770971426421
$ws.Range("E10").Value = 0.00131631144289357
$ws.Range("D11").Value = 0.02953273710448278
$ws.Range("E11").Value = -0.02160007010164755
$ws.Range("D12").Value = 0.01755750063814445
$ws.Range("E12").Value = -0.01112153289561835
$ws.Range("D13").Value = 0.0166535116765731
$ws.Range("E13").Value = -0.003001942433339178
$ws.Range("D14").Value = 0.008194598180729568
$ws.Range("E14").Value = 0.006631832797427739
$ws.Range("D15").Value = 0.007951809213041439
$ws.Range("E15").Value = 0.01301775147928996
$ws.Range("D16").Value = 0.03168254871953522
$ws.Range("E16").Value = 0.03042993985297415
$ws.Range("D17").Value = 0.03190875181297722
$ws.Range("E17").Value = -0.003011837220704638
$ws.Range("D18").Value = 0.03192651399036137
$ws.Range("E18").Value = -0.00498132004981322
$ws.Range("D19").Value = 0.03313857674378118
$ws.Range("E19").Value = -0.001789023065618789
$ws.Range("D20").Value = 0.02664844181004822
$ws.Range("E20").Value = 0.009022529839677995
$ws.Range("D21").Value = 0.03116168169001825
$ws.Range("E21").Value = -0.0009965573473454947
$ws.Range("D22").Value = 0.03328784961199641
$ws.Range("E22").Value = 0.01526571891995054
$ws.Range("D23").Value = 0.03154562703426925
$ws.Range("E23").Value = 0.01133583915040859
$ws.Range("D24").Value = 0.01689924140210629
$ws.Range("E24").Value = 0.001002338790511148
$ws.Range("D25").Value = 0.01570247058948079
$ws.Range("E25").Value = -0.01573151546932372
$ws.Range("D26").Value = 0.03178841600195737
$ws.Range("E26").Value = 0.00451450562462985
$ws.Range("D27").Value = 0.03120238177859387
$ws.Range("E27").Value = 0.00738903256452228
$ws.Range("D28").Value = 0.03276886466750028
$ws.Range("E28").Value = 0.01026653504442243
$ws.Range("D29").Value = 0.0314496406982065
$ws.Range("E29").Value = 0.01324057450628358
$ws.Range("D30").Value = 0.03368226405417581
$ws.Range("E30").Value = 0.005797303904449436
$ws.Range("D31").Value = 0.03185734736584557
$ws.Range("E31").Value = 0.0002104672372666982
$ws.Range("D32").Value = 0.03484586311830905
$ws.Range("E32").Value = -0.005104107591347384
$ws.Range("D33").Value = 0.03004548764234739
$ws.Range("E33").Value = 0.03126969771712029
$ws.Range("D34").Value = 0.04566585227279411
$ws.Range("E34").Value = -0.000347745450330339
$ws.Range("D35").Value = 0.03101958427094496
$ws.Range("E35").Value = 0.007478081485301713
$ws.Range("D36").Value = 0.03350628910472745
$ws.Range("E36").Value = 0.009338444903175036
$ws.Range("E37").Value = 0.005030578000073005

# Restore sheet protection.
$ws.Protect()
